$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph that currently sits right
#    after the "Play Fu Fortunes Megaways Free | Game Review" heading.
# ------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Meta description:") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ------------------------------------------------------------------
# 2. Find the paragraph that currently holds the "Create a feature
#    image..." instructions (italic run) near the end of the document.
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Create a feature image") {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    # Insert a brand-new paragraph right before this one, containing an
    # empty run followed by a bold run with the title text.
    $insertRange = $targetPara.Range.Duplicate
    $insertRange.Collapse(1)
    $insertRange.InsertParagraphBefore()

    $titleRange = $targetPara.Range.Duplicate
    $titleRange.Collapse(1)
    $titleRange.MoveEnd(1, 0)
    $titleRange.Text = "Play Fu Fortunes Megaways Free | Game Review"
    $titleRange.Bold = 1
    $titleRange.Italic = 0

    # Replace the big italic instructions text with the meta description
    # copy, keeping the italic formatting already on that run.
    $d.Content.Find.Execute(
        "Create a feature image for ""Fu Fortunes Megaways"" to add visual appeal to the game review. The image should be in cartoon style and feature a happy Maya warrior with glasses. The warrior should be holding a golden scroll with the game title ""Fu Fortunes Megaways"" written in bold letters. The background should be bright and vibrant, with Asian-inspired elements such as red lanterns, gold coins, and cherry blossom trees. The image should be eye-catching and show the fun and exciting nature of the game.",
        $true, $false, $false, $false, $false, $true, 1, $false,
        "Learn all about Fu Fortunes Megaways casino game and play for free. Read our review on this exciting slot with multiple bonus functions and Megaways mechanic.",
        2
    )
}
